# [PHOENIX-5848] - UI updated create new property approval details
#
# Replace the old single "PTISCommissioner" login-test row (row 83) with
# five PTIS role rows - Junior Assistant, Bill Collector, Revenue
# Inspector, Revenue Officer and Commissioner (the last one is the role
# that used to live alone on row 83, now pushed down to row 87). Each new
# row carries a literal text "0" in column D (replacing the old =FALSE()
# formula) and a "kurnool_eGov@123" mailto hyperlink in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- columns A & B: role name + role code ------------------------------
$ws.Range("A83").Value = "PTISJuniorAssistant"
$ws.Range("B83").Value = "P0001"
$ws.Range("A84").Value = "PTISBillCollector"
$ws.Range("B84").Value = "P0002"
$ws.Range("A85").Value = "PTISRevenueInspector"
$ws.Range("B85").Value = "P0003"
$ws.Range("A86").Value = "PTISRevenueOfficer"
$ws.Range("B86").Value = "P0004"
$ws.Range("A87").Value = "PTISCommissioner"
$ws.Range("B87").Value = "P0005"

$ws.Range("A83:B87").Style = "Normal"

# --- column C: email / hyperlink (rows 84-87 are brand new cells) ------
$ws.Range("C84").Value = "kurnool_eGov@123"
$ws.Range("C85").Value = "kurnool_eGov@123"
$ws.Range("C86").Value = "kurnool_eGov@123"
$ws.Range("C87").Value = "kurnool_eGov@123"

$ws.Hyperlinks.Add($ws.Range("C84"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")
$ws.Hyperlinks.Add($ws.Range("C85"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")
$ws.Hyperlinks.Add($ws.Range("C86"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")
$ws.Hyperlinks.Add($ws.Range("C87"), "mailto:kurnool_eGov@123", "", "", "kurnool_eGov@123")

# Adding a hyperlink auto-applies Excel's "Hyperlink" look (underline +
# theme colour); put the cells back to the plain text style used by the
# rest of column C (same as C83).
$rangeC = $ws.Range("C84:C87")
$rangeC.Font.Name = "Arial"
$rangeC.Font.Size = 11
$rangeC.Font.Underline = $false
$rangeC.Font.Color = 0
$rangeC.NumberFormat = "@"

# --- column D: text "0" flag, right aligned -----------------------------
$rangeD = $ws.Range("D83:D87")
$rangeD.Style = "Excel Built-in Explanatory Text"
$rangeD.Font.Name = "Arial"
$rangeD.Font.Size = 11
$rangeD.NumberFormat = "@"
$rangeD.HorizontalAlignment = -4152

$ws.Range("D83").Value = "0"
$ws.Range("D84").Value = "0"
$ws.Range("D85").Value = "0"
$ws.Range("D86").Value = "0"
$ws.Range("D87").Value = "0"

# --- selection, matching the author's final cursor position ------------
$ws.Range("I83").Select()

Write-Host "done"
